$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the data table. Insert a new
# row at row 15 (shifting the existing rows 15-103 down to 16-104, and
# extending the sheet's used range to A1:R104), then populate it with
# the new record's values.
$ws.Rows.Item(15).EntireRow.Insert()

$ws.Range("A15").Value = 7
$ws.Range("B15").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C15").Value = "Ñuble"
$ws.Range("D15").Value = 44901
$ws.Range("E15").Value = 16
$ws.Range("F15").Value = 100112031
$ws.Range("G15").Value = "Poroto verde"
$ws.Range("H15").Value = "Sin especificar"
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 100
$ws.Range("K15").Value = 30000
$ws.Range("L15").Value = 32000
$ws.Range("M15").Value = 31000
$ws.Range("N15").Value = "$/saco 25 kilos"
$ws.Range("O15").Value = "Región del Maule"
$ws.Range("P15").Value = 1240
$ws.Range("Q15").Value = 25
$ws.Range("R15").Value = "Hortaliza"
